$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("K17").Value = 2627.0001
$ws.Range("H17").Value = 790.4583
$ws.Range("M17").Value = -2459.0001
$ws.Range("I17").Value = 875.6667
$ws.Range("N43").Value = -3138
$ws.Range("I43").Value = 3997.3333
$ws.Range("H43").Value = 3748
$ws.Range("M43").Value = -3928.3333
$ws.Range("J43").Value = 3000
$ws.Range("L43").Value = 3000
$ws.Range("K43").Value = 3997.3333
$ws.Range("H53").Value = 225.22223
$ws.Range("J53").Value = 366.85715
$ws.Range("N53").Value = -1640.85715
$ws.Range("M53").Value = 501.90909
$ws.Range("L53").Value = 366.85715
$ws.Range("K53").Value = 135.09091
$ws.Range("I53").Value = 135.09091
$ws.Range("L86").Value = 6000
$ws.Range("H86").Value = 5642.857
$ws.Range("N86").Value = -8246
$ws.Range("J86").Value = 6000
$ws.Range("M86").Value = -3627
$ws.Range("K86").Value = 4750
$ws.Range("I86").Value = 4750
$ws.Range("M89").Value = -18134
$ws.Range("H89").Value = 5642.857
$ws.Range("K89").Value = 23750
$ws.Range("J89").Value = 6000
$ws.Range("I89").Value = 4750
$ws.Range("N89").Value = -41232
$ws.Range("L89").Value = 30000
$ws.Range("M106").Value = -1862
$ws.Range("I106").Value = 2493
$ws.Range("K106").Value = 2493
$ws.Range("H106").Value = 2469.75
$ws.Range("J111").Value = 5508.4287
$ws.Range("N111").Value = -22659.2861
$ws.Range("H111").Value = 4732
$ws.Range("K111").Value = 6043.5
$ws.Range("L111").Value = 16525.2861
$ws.Range("I111").Value = 2014.5
$ws.Range("M111").Value = -2976.5
$ws.Range("L113").Value = 3349
$ws.Range("H113").Value = 2766.3333
$ws.Range("K113").Value = 2475
$ws.Range("M113").Value = 779
$ws.Range("N113").Value = -9857
$ws.Range("I113").Value = 2475
$ws.Range("J113").Value = 3349
$ws.Range("M129").Value = 211.1818999999996
$ws.Range("H129").Value = 1829.2142
$ws.Range("K129").Value = 4788.8181
$ws.Range("I129").Value = 1596.2727
$ws.Range("H137").Value = 43021360
$ws.Range("K137").Value = 176477790
$ws.Range("N137").Value = -13921603.5
$ws.Range("M137").Value = -176475240
$ws.Range("I137").Value = 58825930
$ws.Range("L137").Value = 13916503.5
$ws.Range("J137").Value = 4638834.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("L96").Value = 49344
$ws.Range("N96").Value = -54836
$ws.Range("H96").Value = 49344
$ws.Range("J96").Value = 49344
$ws.Range("H110").Value = 2278.2144
$ws.Range("I110").Value = 1982.5
$ws.Range("M110").Value = 62.5
$ws.Range("K110").Value = 1982.5
$ws.Range("H122").Value = 3010.087
$ws.Range("I122").Value = 2307.7646
$ws.Range("J122").Value = 5000
$ws.Range("N122").Value = -19900
$ws.Range("K122").Value = 6923.293799999999
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -4473.293799999999
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("H123").Value = 0
$ws.Range("N123").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1339.5834
$ws.Range("N20").Value = -1813
$ws.Range("J20").Value = 1319
$ws.Range("L20").Value = 1319
$ws.Range("H107").Value = 3122.1853
$ws.Range("M107").Value = -772.3890000000001
$ws.Range("K107").Value = 2692.389
$ws.Range("I107").Value = 2692.389
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("K16").Value = 1190.3334
$ws.Range("I16").Value = 1190.3334
$ws.Range("H16").Value = 1345.1666
$ws.Range("M16").Value = -903.3334
$ws.Range("M31").Value = -12307.875
$ws.Range("K31").Value = 12602.875
$ws.Range("J31").Value = 54472.57
$ws.Range("H31").Value = 32142.066
$ws.Range("L31").Value = 54472.57
$ws.Range("I31").Value = 12602.875
$ws.Range("N31").Value = -55062.57
$ws.Range("I34").Value = 12602.875
$ws.Range("H34").Value = 32142.066
$ws.Range("N34").Value = -54876.57
$ws.Range("K34").Value = 12602.875
$ws.Range("J34").Value = 54472.57
$ws.Range("M34").Value = -12400.875
$ws.Range("L34").Value = 54472.57
$ws.Range("H113").Value = 1345.1666
$ws.Range("K113").Value = 1190.3334
$ws.Range("M113").Value = 979.6666
$ws.Range("I113").Value = 1190.3334
$ws.Range("H122").Value = 5373.524
$ws.Range("I122").Value = 2387.3076
$ws.Range("K122").Value = 7161.9228
$ws.Range("M122").Value = -4711.9228
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 17987.8
$ws.Range("N74").Value = -56085.39999999999
$ws.Range("L74").Value = 53963.39999999999
$ws.Range("J74").Value = 17987.8
$ws.Range("H77").Value = 17987.8
$ws.Range("N77").Value = -172498.2
$ws.Range("L77").Value = 161890.2
$ws.Range("J77").Value = 17987.8
$ws.Range("J88").Value = 12617
$ws.Range("N88").Value = -38707
$ws.Range("H88").Value = 12617
$ws.Range("L88").Value = 37851
$ws.Range("L91").Value = 37851
$ws.Range("N91").Value = -40815
$ws.Range("H91").Value = 12617
$ws.Range("J91").Value = 12617
$ws.Range("H94").Value = 4229.357
$ws.Range("J99").Value = 6435.2856
$ws.Range("N99").Value = -23797.8568
$ws.Range("H99").Value = 4176.0835
$ws.Range("L99").Value = 19305.8568
$ws.Range("H120").Value = 26655.166
$ws.Range("I120").Value = 9899.5
$ws.Range("J120").Value = 35033
$ws.Range("K120").Value = 29698.5
$ws.Range("N120").Value = -114775
$ws.Range("L120").Value = 105099
$ws.Range("M120").Value = -24860.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H75").Value = 57499
$ws.Range("J75").Value = 57499
$ws.Range("N75").Value = -59247
$ws.Range("L75").Value = 57499
$ws.Range("L78").Value = 172497
$ws.Range("J78").Value = 57499
$ws.Range("N78").Value = -181233
$ws.Range("H78").Value = 57499
$ws.Range("N80").Value = -7000
$ws.Range("L80").Value = 5004
$ws.Range("M80").Value = -515484.6
$ws.Range("J80").Value = 5004
$ws.Range("H80").Value = 398449.06
$ws.Range("I80").Value = 516482.6
$ws.Range("K80").Value = 516482.6
$ws.Range("K83").Value = 2582413
$ws.Range("N83").Value = -35004
$ws.Range("I83").Value = 516482.6
$ws.Range("H83").Value = 398449.06
$ws.Range("J83").Value = 5004
$ws.Range("L83").Value = 25020
$ws.Range("M83").Value = -2577421
$ws.Range("H102").Value = 2969.8372
$ws.Range("I102").Value = 2053.6
$ws.Range("M102").Value = -431.5999999999999
$ws.Range("K102").Value = 2053.6
$ws.Range("H107").Value = 23794.8
$ws.Range("M107").Value = -37305.223
$ws.Range("K107").Value = 39225.223
$ws.Range("N107").Value = -4489.1667
$ws.Range("I107").Value = 39225.223
$ws.Range("L107").Value = 649.1667
$ws.Range("J107").Value = 649.1667
$ws.Range("H122").Value = 6530
$ws.Range("I122").Value = 3920.2
$ws.Range("K122").Value = 11760.6
$ws.Range("M122").Value = -9310.599999999999
$ws.Range("J123").Value = 59999.5
$ws.Range("L123").Value = 59999.5
$ws.Range("H123").Value = 59999.5
$ws.Range("N123").Value = -64899.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M22").Value = -251.2
$ws.Range("I22").Value = 546.2
$ws.Range("K22").Value = 546.2
$ws.Range("H22").Value = 1022.625
$ws.Range("K27").Value = 546.2
$ws.Range("H27").Value = 1022.625
$ws.Range("I27").Value = 546.2
$ws.Range("M27").Value = -439.2
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("L46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("H46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("I93").Value = 1115.8462
$ws.Range("H93").Value = 1261.28
$ws.Range("K93").Value = 1115.8462
$ws.Range("M93").Value = 132.1538
$ws.Range("J136").Value = 4583.5
$ws.Range("H136").Value = 3066.3103
$ws.Range("N136").Value = -18850.5
$ws.Range("K136").Value = 6803.3688
$ws.Range("M136").Value = -4253.3688
$ws.Range("I136").Value = 2267.7896
$ws.Range("L136").Value = 13750.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3915.5
$ws.Range("L81").Value = 12332
$ws.Range("J81").Value = 6166
$ws.Range("N81").Value = -14454
$ws.Range("H84").Value = 3915.5
$ws.Range("L84").Value = 61660
$ws.Range("J84").Value = 6166
$ws.Range("N84").Value = -72268
$ws.Range("L114").Value = 90398
$ws.Range("N114").Value = -99076
$ws.Range("J114").Value = 90398
$ws.Range("H114").Value = 90398
$ws.Range("H132").Value = 11805130
$ws.Range("N132").Value = -26434.5005
$ws.Range("I132").Value = 14333274
$ws.Range("J132").Value = 7124.8335
$ws.Range("K132").Value = 42999822
$ws.Range("M132").Value = -42997292
$ws.Range("L132").Value = 21374.5005
$ws.Range("J136").Value = 6390.6665
$ws.Range("H136").Value = 12830745
$ws.Range("N136").Value = -24271.9995
$ws.Range("K136").Value = 50034153
$ws.Range("M136").Value = -50031603
$ws.Range("I136").Value = 16678051
$ws.Range("L136").Value = 19171.9995
